# "Todos los botones funcionales"
# The sheet had 5 data rows (rows 2-6). Rows 3, 5 and 6 contained leftover
# test/junk entries that are removed, leaving only the data rows 2 and 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(3).ClearContents()
$ws.Rows(5).ClearContents()
$ws.Rows(6).ClearContents()
